# Auto-upload VRF Excel file: append a new "qwe" sheet with the standard
# Outdoor/Indoor model header row used by every other sheet in this workbook.

$wb = $excel.ActiveWorkbook

# Remember the sheet that was active before we start, so the workbook's
# active-tab selection is unchanged by adding/naming the new sheet.
$originalActive = $wb.ActiveSheet

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab strip (matches the <sheets> ordering in workbook.xml).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "qwe"

# Match the page-margin template (inches) used by all the other sheets.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Standard header row: Outdoor/Indoor model, quantity, serial(s).
$headers = @(
    "Outdoor Model",
    "Outdoor Quantity",
    "Outdoor Serial(s)",
    "Indoor Model",
    "Indoor Quantity",
    "Indoor Serial(s)"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

$ws.Range("A1").Select()

# Restore the original active sheet/tab so only the new sheet is added.
$originalActive.Activate()
